$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the coin snapshot table (price + 1h volume change, and a few
# ranking swaps) to the latest scrape. Every cell in columns B:E is stored
# as text in the source data, so each write is prefixed with a leading
# apostrophe to stop Excel from auto-coercing number-looking strings (e.g.
# "1.00", "0.650") into actual numbers that would drop trailing zeros; the
# explicit Style reset afterwards clears the transient quote-prefix glyph
# style so the cell formatting stays identical to every other data cell.

$ws.Range("D2").Value = "'37.459.11"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +5.60%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.043.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +3.22%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.04%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'252.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +4.88%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.650"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +3.15%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'65.40"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +18.72%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.00%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +6.84%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'59.11"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.23%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +4.95%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +1.03%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.908"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +2.90%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'15.14"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +7.94%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'2.336.43"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +3.28%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "'Avalanche"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'20.92"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +24.40%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "'Polkadot"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'5.60"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +8.33%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.038.94"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +3.27%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'37.374.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +5.76%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'73.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +4.89%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.0₃0876"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +6.00%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +8.56%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'236.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +2.68%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +23.42%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -0.19%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +1.89%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'9.57"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +6.34%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'165.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +1.87%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'19.88"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +3.18%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.122"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +2.89%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'5.22"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +10.13%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'1.22"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +8.84%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.112"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +27.09%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'4.73"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +12.43%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.0617"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +5.98%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'2.46"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +13.51%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.06%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +0.56%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'5.97"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +24.28%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +20.13%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +5.66%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +4.29%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'RenderToken"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'2.74"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +23.85%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "'VeChain"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'0.0219"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +6.31%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  +7.09%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "'InjectiveProtocol"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'17.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +12.92%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "'FraxShare"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'8.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +11.14%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'95.61"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +7.62%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.419.38"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +3.48%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  +2.31%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'47.47"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +6.90%  "
$ws.Range("E51").Style = "Normal"
